# Auto-generated Excel COM-interop script to apply diff changes to Atomos_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1424.8064
$ws.Range("I8").Value = 214.75
$ws.Range("K8").Value = 644.25
$ws.Range("M8").Value = -505.25
$ws.Range("H19").Value = 10204425
$ws.Range("I19").Value = 20408416
$ws.Range("J19").Value = 434.2857
$ws.Range("K19").Value = 20408416
$ws.Range("L19").Value = 434.2857
$ws.Range("M19").Value = -20408241
$ws.Range("N19").Value = -784.2857
$ws.Range("H32").Value = 679.8
$ws.Range("I32").Value = 633
$ws.Range("J32").Value = 750
$ws.Range("K32").Value = 633
$ws.Range("L32").Value = 750
$ws.Range("M32").Value = -307
$ws.Range("N32").Value = -1402
$ws.Range("H53").Value = 359.72415
$ws.Range("I53").Value = 434.36365
$ws.Range("K53").Value = 434.36365
$ws.Range("M53").Value = 202.63635
$ws.Range("H87").Value = 22773.115
$ws.Range("J87").Value = 22773.115
$ws.Range("L87").Value = 22773.115
$ws.Range("N87").Value = -25269.115
$ws.Range("H90").Value = 22773.115
$ws.Range("J90").Value = 22773.115
$ws.Range("L90").Value = 68319.345
$ws.Range("N90").Value = -80799.345
$ws.Range("H116").Value = 4901.5293
$ws.Range("I116").Value = 4408
$ws.Range("J116").Value = 5606.5713
$ws.Range("K116").Value = 4408
$ws.Range("L116").Value = 5606.5713
$ws.Range("M116").Value = -966
$ws.Range("N116").Value = -12490.5713
$ws.Range("H121").Value = 843.5
$ws.Range("J121").Value = 1137
$ws.Range("L121").Value = 3411
$ws.Range("N121").Value = -6905
$ws.Range("H129").Value = 4238291.5
$ws.Range("I129").Value = 35715372
$ws.Range("J129").Value = 992.2308
$ws.Range("K129").Value = 107146116
$ws.Range("L129").Value = 2976.6924
$ws.Range("M129").Value = -107141116
$ws.Range("N129").Value = -12976.6924
$ws.Range("H135").Value = 2433.3333
$ws.Range("I135").Value = 3475
$ws.Range("K135").Value = 31275
$ws.Range("M135").Value = -28740
$ws.Range("H138").Value = 3598.0645
$ws.Range("I138").Value = 4045.5642
$ws.Range("J138").Value = 2839.261
$ws.Range("K138").Value = 12136.6926
$ws.Range("L138").Value = 8517.782999999999
$ws.Range("M138").Value = -6996.692599999998
$ws.Range("N138").Value = -18797.783
$ws.Range("H141").Value = 620106.0600000001
$ws.Range("I141").Value = 2747
$ws.Range("J141").Value = 844600.25
$ws.Range("K141").Value = 8241
$ws.Range("L141").Value = 2533800.75
$ws.Range("M141").Value = -3061
$ws.Range("N141").Value = -2544160.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2549.15
$ws.Range("I63").Value = 1936.5
$ws.Range("J63").Value = 4999.75
$ws.Range("K63").Value = 1936.5
$ws.Range("L63").Value = 4999.75
$ws.Range("M63").Value = -1250.5
$ws.Range("N63").Value = -6371.75
$ws.Range("H66").Value = 2549.15
$ws.Range("I66").Value = 1936.5
$ws.Range("J66").Value = 4999.75
$ws.Range("K66").Value = 9682.5
$ws.Range("L66").Value = 24998.75
$ws.Range("M66").Value = -6250.5
$ws.Range("N66").Value = -31862.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1391.3846
$ws.Range("I20").Value = 1248.8889
$ws.Range("J20").Value = 1712
$ws.Range("K20").Value = 1248.8889
$ws.Range("L20").Value = 1712
$ws.Range("M20").Value = -1001.8889
$ws.Range("N20").Value = -2206
$ws.Range("H80").Value = 561.5714
$ws.Range("I80").Value = 686
$ws.Range("K80").Value = 686
$ws.Range("M80").Value = 312
$ws.Range("H83").Value = 561.5714
$ws.Range("I83").Value = 686
$ws.Range("K83").Value = 3430
$ws.Range("M83").Value = 1562
$ws.Range("H134").Value = 2569.0588
$ws.Range("I134").Value = 1777.8334
$ws.Range("K134").Value = 5333.5002
$ws.Range("M134").Value = -2798.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 23857.857
$ws.Range("I4").Value = 4250.25
$ws.Range("J4").Value = 50001.332
$ws.Range("K4").Value = 4250.25
$ws.Range("L4").Value = 50001.332
$ws.Range("M4").Value = -4138.25
$ws.Range("N4").Value = -50225.332
$ws.Range("H19").Value = 1109.1666
$ws.Range("I19").Value = 331
$ws.Range("K19").Value = 331
$ws.Range("M19").Value = -161
$ws.Range("H24").Value = 1109.1666
$ws.Range("I24").Value = 331
$ws.Range("K24").Value = 331
$ws.Range("M24").Value = -161
$ws.Range("H31").Value = 2129950
$ws.Range("I31").Value = 2633094.5
$ws.Range("J31").Value = 5562.5557
$ws.Range("K31").Value = 2633094.5
$ws.Range("L31").Value = 5562.5557
$ws.Range("M31").Value = -2632799.5
$ws.Range("N31").Value = -6152.5557
$ws.Range("H34").Value = 2129950
$ws.Range("I34").Value = 2633094.5
$ws.Range("J34").Value = 5562.5557
$ws.Range("K34").Value = 2633094.5
$ws.Range("L34").Value = 5562.5557
$ws.Range("M34").Value = -2632892.5
$ws.Range("N34").Value = -5966.5557
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1134.3158
$ws.Range("I122").Value = 526.5
$ws.Range("J122").Value = 1809.6666
$ws.Range("K122").Value = 4738.5
$ws.Range("L122").Value = 16286.9994
$ws.Range("M122").Value = -2288.5
$ws.Range("N122").Value = -21186.9994
$ws.Range("H131").Value = 1530.8096
$ws.Range("J131").Value = 1204.8125
$ws.Range("L131").Value = 3614.4375
$ws.Range("N131").Value = -13694.4375
$ws.Range("H132").Value = 2045.3846
$ws.Range("I132").Value = 1714.1428
$ws.Range("J132").Value = 2431.8333
$ws.Range("K132").Value = 15427.2852
$ws.Range("L132").Value = 21886.4997
$ws.Range("M132").Value = -12897.2852
$ws.Range("N132").Value = -26946.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 235005
$ws.Range("I18").Value = 1000000
$ws.Range("J18").Value = 82006
$ws.Range("K18").Value = 1000000
$ws.Range("L18").Value = 82006
$ws.Range("M18").Value = -999707
$ws.Range("N18").Value = -82592
$ws.Range("H20").Value = 44003
$ws.Range("J20").Value = 44003
$ws.Range("L20").Value = 44003
$ws.Range("N20").Value = -44493
$ws.Range("H21").Value = 80007
$ws.Range("J21").Value = 80007
$ws.Range("L21").Value = 80007
$ws.Range("N21").Value = -80353
$ws.Range("H22").Value = 93342.336
$ws.Range("J22").Value = 93342.336
$ws.Range("L22").Value = 93342.336
$ws.Range("N22").Value = -94400.336
$ws.Range("H24").Value = 19252.125
$ws.Range("J24").Value = 19252.125
$ws.Range("L24").Value = 19252.125
$ws.Range("N24").Value = -19598.125
$ws.Range("H30").Value = 80007
$ws.Range("J30").Value = 80007
$ws.Range("L30").Value = 80007
$ws.Range("N30").Value = -80217
$ws.Range("H113").Value = 2215.7144
$ws.Range("I113").Value = 1627.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1627.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 542.5
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2686.9033
$ws.Range("J132").Value = 3151.9
$ws.Range("L132").Value = 9455.700000000001
$ws.Range("N132").Value = -14515.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 57679.332
$ws.Range("J31").Value = 57679.332
$ws.Range("L31").Value = 57679.332
$ws.Range("N31").Value = -58375.332
$ws.Range("H132").Value = 142049.6
$ws.Range("I132").Value = 176091.89
$ws.Range("J132").Value = 3448.7856
$ws.Range("K132").Value = 528275.67
$ws.Range("L132").Value = 10346.3568
$ws.Range("M132").Value = -525745.67
$ws.Range("N132").Value = -15406.3568
$ws.Range("H136").Value = 1416.7333
$ws.Range("I136").Value = 674.381
$ws.Range("K136").Value = 2023.143
$ws.Range("M136").Value = 526.857
